$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.002833900968242052
$ws.Range("C2").Value = -5864.14666073515
$ws.Range("D2").Value = -115.1647534024966
$ws.Range("E2").Value = 51410.93713904079
$ws.Range("F2").Value = 1.596731641488282
$ws.Range("G2").Value = -1.233552899527711
$ws.Range("H2").Value = 0.06949803370384533
$ws.Range("I2").Value = -1.68783007370961
$ws.Range("J2").Value = 4.708258456577482
$ws.Range("K2").Value = 288
$ws.Range("L2").Value = 160.13036782199
$ws.Range("N2").Value = 4.708258456577482
$ws.Range("O2").Value = 5.374925123244149
$ws.Range("B3").Value = -346.4544766790139
$ws.Range("C3").Value = -5138.16875532128
$ws.Range("D3").Value = 13318.72736983711
$ws.Range("E3").Value = 0.0001857736446508796
$ws.Range("F3").Value = -0.6405017168990723
$ws.Range("G3").Value = -1.551981812672806
$ws.Range("H3").Value = -1.341692803832187
$ws.Range("I3").Value = 1.989389935750627
$ws.Range("J3").Value = 4.708273077746254
$ws.Range("K3").Value = 581
$ws.Range("L3").Value = 0.39099350149622
$ws.Range("N3").Value = 4.708273077746261
$ws.Range("O3").Value = 5.374939744412928
$ws.Range("B4").Value = -12937.70730302205
$ws.Range("C4").Value = 0.0009968436357354659
$ws.Range("D4").Value = 13315.72371925506
$ws.Range("E4").Value = 22073.55800362063
$ws.Range("F4").Value = -1.017822263739436
$ws.Range("G4").Value = 1.729893523390114
$ws.Range("H4").Value = -1.013779453922141
$ws.Range("I4").Value = -1.695854597040008
$ws.Range("J4").Value = 4.708281036686117
$ws.Range("K4").Value = 52
$ws.Range("L4").Value = -10.03184155721563
$ws.Range("N4").Value = 4.708281036718855
$ws.Range("O4").Value = 5.374947703385522
$ws.Range("B5").Value = 0.0007700566648021041
$ws.Range("C5").Value = -1527.798229141081
$ws.Range("D5").Value = -1605.011723441246
$ws.Range("E5").Value = 9174.3002126514
$ws.Range("F5").Value = 1.7644558855894
$ws.Range("G5").Value = -0.8571358208756725
$ws.Range("H5").Value = -1.144144103285236
$ws.Range("I5").Value = -1.140364276522475
$ws.Range("J5").Value = 4.708292990727894
$ws.Range("K5").Value = 668
$ws.Range("L5").Value = -5.327749651820979
$ws.Range("N5").Value = 4.70829299072791
$ws.Range("O5").Value = 5.374959657394577
$ws.Range("B6").Value = 0.001852846804108962
$ws.Range("C6").Value = -0.4202368598012528
$ws.Range("D6").Value = -2522.329220115455
$ws.Range("E6").Value = 12197.20793857272
$ws.Range("F6").Value = 1.635660178183984
$ws.Range("G6").Value = 0.3736159243433201
$ws.Range("H6").Value = -1.720160387742061
$ws.Range("I6").Value = -1.455457572209112
$ws.Range("J6").Value = 4.708295376414575
$ws.Range("K6").Value = 107
$ws.Range("L6").Value = -7.59671641329339
$ws.Range("N6").Value = 4.708295376432961
$ws.Range("O6").Value = 5.374962043099628
$ws.Range("B7").Value = -3497.990009461215
$ws.Range("C7").Value = 0.0006923256595945056
$ws.Range("D7").Value = -3996.847146262389
$ws.Range("E7").Value = 31202.40295610447
$ws.Range("F7").Value = -1.460960037463827
$ws.Range("G7").Value = 1.780521669734684
$ws.Range("H7").Value = -1.750540777755884
$ws.Range("I7").Value = -1.621809611138783
$ws.Range("J7").Value = 4.708307027130176
$ws.Range("K7").Value = 243
$ws.Range("L7").Value = -7.246890860550308
$ws.Range("N7").Value = 4.708307027176368
$ws.Range("O7").Value = 5.374973693843035
$ws.Range("B8").Value = 0.0005149403881635328
$ws.Range("C8").Value = -1103.912773063059
$ws.Range("D8").Value = -177.7457372487873
$ws.Range("E8").Value = 3514.499383008041
$ws.Range("F8").Value = 1.810488400907328
$ws.Range("G8").Value = -1.389026765483945
$ws.Range("H8").Value = -0.1917547081398894
$ws.Range("I8").Value = -0.9945498152204615
$ws.Range("J8").Value = 4.708315618662134
$ws.Range("K8").Value = 730
$ws.Range("L8").Value = 44.71980227531387
$ws.Range("N8").Value = 4.708315618662155
$ws.Range("O8").Value = 5.374982285328822
$ws.Range("B9").Value = 0.0006533019421955117
$ws.Range("C9").Value = 515.8497165115939
$ws.Range("D9").Value = 2500.740130313069
$ws.Range("E9").Value = 70.71432700333804
$ws.Range("F9").Value = 1.768126010259178
$ws.Range("G9").Value = -1.34773997977639
$ws.Range("H9").Value = -0.9792683194705409
$ws.Range("I9").Value = 0.08659808241459199
$ws.Range("J9").Value = 4.708322286540437
$ws.Range("K9").Value = 868
$ws.Range("L9").Value = -128.5995452926111
$ws.Range("N9").Value = 4.708322738314763
$ws.Range("O9").Value = 5.37498940498143
$ws.Range("B10").Value = -119.2465728251387
$ws.Range("C10").Value = -507.4873829367511
$ws.Range("D10").Value = 0.0001347548805125384
$ws.Range("E10").Value = 1770.022216781039
$ws.Range("F10").Value = -1.502315923047965
$ws.Range("G10").Value = -0.07101887494854675
$ws.Range("H10").Value = 1.997798873210352
$ws.Range("I10").Value = -0.735354022113498
$ws.Range("J10").Value = 4.708327547191402
$ws.Range("K10").Value = 606
$ws.Range("L10").Value = 312.3899551840778
$ws.Range("N10").Value = 4.708327547191459
$ws.Range("O10").Value = 5.374994213858126
$ws.Range("B11").Value = -277.5314346699907
$ws.Range("C11").Value = 0.003968030574984961
$ws.Range("D11").Value = -6.570812047035906
$ws.Range("E11").Value = 2873.465122083891
$ws.Range("F11").Value = -0.3157809265194105
$ws.Range("G11").Value = 1.515354533756579
$ws.Range("H11").Value = 0.2601848759091361
$ws.Range("I11").Value = -0.9331410691209625
$ws.Range("J11").Value = 4.708329508185642
$ws.Range("K11").Value = 477
$ws.Range("L11").Value = 50.81299256860409
$ws.Range("N11").Value = 4.708329508185649
$ws.Range("O11").Value = 5.374996174852316
